$d = $word.ActiveDocument

# 1. Title / heading text (appears twice: Heading1 and bold run near the end)
$d.Content.Find.Execute(
    "Play Arcane Gems for Free- Review of this High-Rewarding Slot",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Arcane Gems Free: Review of Symbol Locking & High Rewards",
    2)

# Find.Execute only updates the first match by default via a single call
# on Content; run it again to catch the second occurrence (bold run).
$d.Content.Find.Execute(
    "Play Arcane Gems for Free- Review of this High-Rewarding Slot",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Arcane Gems Free: Review of Symbol Locking & High Rewards",
    2)

# 2. "What we don't like" bullet list items
$d.Content.Find.Execute(
    "High Volatility",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Volatility and RTP",
    2)

$d.Content.Find.Execute(
    "Low Base Game Wins",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paytable and Payouts",
    2)

# 3. Meta description (italic run)
$d.Content.Find.Execute(
    "Read our review of Arcane Gems and play for free. Find out about its high rewards, symbol locking feature, and respin feature. Available for free play.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Arcane Gems for free and enjoy the symbol locking feature and high rewards.",
    2)
